# Update countries & provincias Spain
# Refresh COVID stats snapshot (4 May 2020, 12:08 -> 12:33) and the
# resulting reorder of two country pairs whose rank changed because of it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Timestamp banner (row 1)
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 4 de Mayo de 2020 a las 12:33"

# --- Plain numeric refresh (country stays on the same row) ---

# Row 4: Estados Unidos
$ws.Cells.Item(4, 2).Value = 1188826
$ws.Cells.Item(4, 3).Value = 704
$ws.Cells.Item(4, 5).Value = 941626
$ws.Cells.Item(4, 7).Value = 8
$ws.Cells.Item(4, 8).Value = 68606

# Row 20: Suiza
$ws.Cells.Item(20, 2).Value = 29981
$ws.Cells.Item(20, 3).Value = 76
$ws.Cells.Item(20, 5).Value = 3719

# Row 31: Israel
$ws.Cells.Item(31, 6).Value = 93

# Row 50: Australia
$ws.Cells.Item(50, 2).Value = 6825
$ws.Cells.Item(50, 3).Value = 24
$ws.Cells.Item(50, 4).Value = 5859
$ws.Cells.Item(50, 5).Value = 871
$ws.Cells.Item(50, 6).Value = 28

# Row 76: Bosnia y Herzegovina
$ws.Cells.Item(76, 2).Value = 1926
$ws.Cells.Item(76, 3).Value = 69
$ws.Cells.Item(76, 4).Value = 855
$ws.Cells.Item(76, 5).Value = 993
$ws.Cells.Item(76, 7).Value = 1
$ws.Cells.Item(76, 8).Value = 78

# Row 85: Eslovenia
$ws.Cells.Item(85, 5).Value = 1101
$ws.Cells.Item(85, 6).Value = 20
$ws.Cells.Item(85, 7).Value = 1
$ws.Cells.Item(85, 8).Value = 97

# --- Rows whose ranking swapped with the row below/above ---

# Row 55/56 pair (Kuwait / Marruecos): Marruecos's refreshed total (5000)
# overtakes Kuwait's unchanged total (4983), so Marruecos now sorts above
# Kuwait; Kuwait's figures are carried down to row 56 unchanged.
$ws.Cells.Item(55, 1).Value = "Marruecos"
$ws.Cells.Item(55, 2).Value = 5000
$ws.Cells.Item(55, 3).Value = 97
$ws.Cells.Item(55, 4).Value = 1565
$ws.Cells.Item(55, 5).Value = 3258
$ws.Cells.Item(55, 6).Value = 1
$ws.Cells.Item(55, 7).Value = 3
$ws.Cells.Item(55, 8).Value = 177

$ws.Cells.Item(56, 1).Value = "Kuwait"
$ws.Cells.Item(56, 2).Value = 4983
$ws.Cells.Item(56, 3).Value = 0
$ws.Cells.Item(56, 4).Value = 1776
$ws.Cells.Item(56, 5).Value = 3169
$ws.Cells.Item(56, 6).Value = 72
$ws.Cells.Item(56, 7).Value = 0
$ws.Cells.Item(56, 8).Value = 38

# Row 100/101 pair (Costa Rica / Libano): Libano's refreshed total (740)
# overtakes Costa Rica's unchanged total (739), so Libano now sorts above
# Costa Rica; Costa Rica's figures are carried down to row 101 unchanged.
$ws.Cells.Item(100, 1).Value = "Libano"
$ws.Cells.Item(100, 2).Value = 740
$ws.Cells.Item(100, 3).Value = 3
$ws.Cells.Item(100, 4).Value = 200
$ws.Cells.Item(100, 5).Value = 515
$ws.Cells.Item(100, 6).Value = 43
$ws.Cells.Item(100, 7).Value = 0
$ws.Cells.Item(100, 8).Value = 25

$ws.Cells.Item(101, 1).Value = "Costa Rica"
$ws.Cells.Item(101, 2).Value = 739
$ws.Cells.Item(101, 3).Value = 0
$ws.Cells.Item(101, 4).Value = 386
$ws.Cells.Item(101, 5).Value = 347
$ws.Cells.Item(101, 6).Value = 6
$ws.Cells.Item(101, 7).Value = 0
$ws.Cells.Item(101, 8).Value = 6

Write-Output "Applied country/provincia refresh"
